# Update cryptocurrency price/volume data (refresh run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to stay a text value (no auto-conversion to
    # number/date) and then drop back to the default "Normal" style
    # so we don't leave a stray number-format override behind.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '29.386.27'
$ws.Range("E2").Value = '  +0.27%  '

Set-TextValue $ws.Range("D3") '1.867.97'
$ws.Range("E3").Value = '  -0.38%  '

$ws.Range("E4").Value = '  +0.00%  '

Set-TextValue $ws.Range("D5") '243.79'
$ws.Range("E5").Value = '  +0.60%  '

Set-TextValue $ws.Range("D6") '0.7040'
$ws.Range("E6").Value = '  -1.16%  '

Set-TextValue $ws.Range("D8") '0.07931'
$ws.Range("E8").Value = '  -1.30%  '

Set-TextValue $ws.Range("D9") '0.3134'
$ws.Range("E9").Value = '  -0.42%  '

Set-TextValue $ws.Range("D10") '24.45'
$ws.Range("E10").Value = '  -2.12%  '

Set-TextValue $ws.Range("D11") '0.07843'
$ws.Range("E11").Value = '  -4.57%  '

Set-TextValue $ws.Range("D12") '1.864.59'
$ws.Range("E12").Value = '  -0.84%  '

Set-TextValue $ws.Range("D13") '93.89'
$ws.Range("E13").Value = '  -1.04%  '

Set-TextValue $ws.Range("D14") '5.180'
$ws.Range("E14").Value = '  -1.30%  '

Set-TextValue $ws.Range("D15") '0.7009'
$ws.Range("E15").Value = '  -1.55%  '

Set-TextValue $ws.Range("D16") '6.529'
$ws.Range("E16").Value = '  +1.93%  '

Set-TextValue $ws.Range("D17") '0.000008399'
$ws.Range("E17").Value = '  -1.40%  '

Set-TextValue $ws.Range("D18") '29.367.96'
$ws.Range("E18").Value = '  +0.20%  '

Set-TextValue $ws.Range("D19") '253.96'

Set-TextValue $ws.Range("D20") '2.113.99'
$ws.Range("E20").Value = '  -0.71%  '

Set-TextValue $ws.Range("D21") '13.11'
$ws.Range("E21").Value = '  -0.93%  '

Set-TextValue $ws.Range("D22") '0.9999'
$ws.Range("E22").Value = '  -0.07%  '

Set-TextValue $ws.Range("D23") '7.647'
$ws.Range("E23").Value = '  -1.50%  '

$ws.Range("E24").Value = '  -0.04%  '

Set-TextValue $ws.Range("D25") '0.1553'
$ws.Range("E25").Value = '  -0.37%  '

Set-TextValue $ws.Range("D26") '9.006'
$ws.Range("E26").Value = '  -0.33%  '

Set-TextValue $ws.Range("D27") '161.66'
$ws.Range("E27").Value = '  -0.38%  '

Set-TextValue $ws.Range("D28") '18.84'
$ws.Range("E28").Value = '  +1.71%  '

Set-TextValue $ws.Range("D29") '1.503'
$ws.Range("E29").Value = '  +0.09%  '

Set-TextValue $ws.Range("D30") '4.316'
$ws.Range("E30").Value = '  -1.95%  '

Set-TextValue $ws.Range("D31") '4.248'
$ws.Range("E31").Value = '  -1.17%  '

Set-TextValue $ws.Range("D32") '1.216'
$ws.Range("E32").Value = '  +3.29%  '

Set-TextValue $ws.Range("D33") '0.05280'
$ws.Range("E33").Value = '  -1.64%  '

Set-TextValue $ws.Range("D34") '1.898'
$ws.Range("E34").Value = '  -2.06%  '

Set-TextValue $ws.Range("D35") '1.175'
$ws.Range("E35").Value = '  -0.16%  '

Set-TextValue $ws.Range("D36") '0.7469'
$ws.Range("E36").Value = '  -2.21%  '

Set-TextValue $ws.Range("D37") '2.711'
$ws.Range("E37").Value = '  +0.81%  '

Set-TextValue $ws.Range("D38") '0.01880'
$ws.Range("E38").Value = '  +0.27%  '

Set-TextValue $ws.Range("D39") '1.275.57'
$ws.Range("E39").Value = '  +0.82%  '

Set-TextValue $ws.Range("D40") '2.766'
$ws.Range("E40").Value = '  +0.50%  '

Set-TextValue $ws.Range("D41") '0.8915'
$ws.Range("E41").Value = '  -2.59%  '

Set-TextValue $ws.Range("D42") '5.996'
$ws.Range("E42").Value = '  -7.07%  '

Set-TextValue $ws.Range("D43") '108.48'
$ws.Range("E43").Value = '  -3.45%  '

Set-TextValue $ws.Range("D44") '71.07'
$ws.Range("E44").Value = '  -3.86%  '

$ws.Range("E45").Value = '  +0.04%  '

Set-TextValue $ws.Range("D46") '0.00000000128'
$ws.Range("E46").Value = '  -4.85%  '

Set-TextValue $ws.Range("D47") '2.015.36'
$ws.Range("E47").Value = '  -0.54%  '

Set-TextValue $ws.Range("D48") '9.599'
$ws.Range("E48").Value = '  +1.28%  '

Set-TextValue $ws.Range("D49") '1.797'
$ws.Range("E49").Value = '  -0.04%  '

$ws.Range("E50").Value = '  -0.84%  '

Set-TextValue $ws.Range("D51") '0.4301'
$ws.Range("E51").Value = '  -1.12%  '
